# Applies the Coinranking crypto price/volume update described by the commit
# 'Updated cryptos list on Wed Oct  4 21:19:44 UTC 2023 with GitHub Actions'.
# Rows 2-51 (Price column D / Volume(1h) column E) get refreshed figures; rows
# 28 and 29 additionally swap places (BinanceUSD <-> EthereumClassic).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.728.36"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "1.646.19"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "213.38"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  +3.58%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "23.13"
$ws.Range("E8").Value = "  -2.03%  "
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "0.0613"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("D12").Value = "1.877.43"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "1.642.74"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "0.563"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").Value = "64.28"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("D17").Value = "27.702.32"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").Value = "230.95"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").Value = "0.0₃0726"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").Value = "7.67"
$ws.Range("E20").Value = "  +2.49%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").Value = "10.04"
$ws.Range("E23").Value = "  +7.17%  "
$ws.Range("E24").Value = "  -3.37%  "
$ws.Range("D25").Value = "149.12"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").Value = "6.98"
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "15.69"
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("D31").Value = "0.0485"
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").Value = "1.445.76"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").Value = "0.571"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  -2.55%  "
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("D40").Value = "0.894"
$ws.Range("E40").Value = "  +13.25%  "
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "5.65"
$ws.Range("E43").Value = "  +2.39%  "
$ws.Range("D44").Value = "2.46"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").Value = "65.65"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D47").Value = "1.787.60"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "86.01"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").Value = "0.0₆0108"
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0990"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -2.09%  "
